$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells are treated as text so numeric-looking
# strings (e.g. "24.17") are not auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.321.96'
$ws.Range("E2").Value = '  +4.04%  '

$ws.Range("D3").Value = '1.735.14'
$ws.Range("E3").Value = '  +3.01%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '220.03'
$ws.Range("E5").Value = '  +1.96%  '

$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '24.17'
$ws.Range("E8").Value = '  +11.50%  '

$ws.Range("E9").Value = '  +4.70%  '

$ws.Range("D10").Value = '0.0638'
$ws.Range("E10").Value = '  +2.19%  '

$ws.Range("D11").Value = '0.0897'
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").Value = '1.978.10'
$ws.Range("E12").Value = '  +3.00%  '

$ws.Range("D13").Value = '1.736.30'
$ws.Range("E13").Value = '  +3.16%  '

$ws.Range("E14").Value = '  +2.85%  '

$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("D16").Value = '67.71'
$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("D17").Value = '28.299.88'
$ws.Range("E17").Value = '  +4.04%  '

$ws.Range("D18").Value = '242.38'
$ws.Range("E18").Value = '  +1.19%  '

$ws.Range("D19").Value = '0.0₃0758'
$ws.Range("E19").Value = '  +1.80%  '

$ws.Range("D20").Value = '8.00'
$ws.Range("E20").Value = '  -1.20%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '4.67'
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("D23").Value = '9.78'
$ws.Range("E23").Value = '  +2.67%  '

$ws.Range("E24").Value = '  +0.15%  '

$ws.Range("D25").Value = '149.81'
$ws.Range("E25").Value = '  +1.03%  '

$ws.Range("D26").Value = '7.56'
$ws.Range("E26").Value = '  +3.80%  '

$ws.Range("E27").Value = '  +1.89%  '

$ws.Range("E28").Value = '  +0.99%  '

$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("E30").Value = '  +3.03%  '

$ws.Range("E31").Value = '  +2.74%  '

$ws.Range("E32").Value = '  +2.01%  '

$ws.Range("D33").Value = '1.507.79'
$ws.Range("E33").Value = '  -4.41%  '

$ws.Range("D34").Value = '3.28'
$ws.Range("E34").Value = '  +1.25%  '

$ws.Range("E35").Value = '  -2.23%  '

$ws.Range("D36").Value = '0.968'
$ws.Range("E36").Value = '  +2.49%  '

$ws.Range("D37").Value = '0.605'
$ws.Range("E37").Value = '  +0.47%  '

$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("E39").Value = '  +1.91%  '

$ws.Range("E40").Value = '  +1.53%  '

$ws.Range("D41").Value = '70.69'
$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("D42").Value = '5.71'
$ws.Range("E42").Value = '  +1.49%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("E44").Value = '  +1.91%  '

$ws.Range("D45").Value = '1.883.01'
$ws.Range("E45").Value = '  +2.79%  '

$ws.Range("E46").Value = '  +1.94%  '

$ws.Range("E47").Value = '  +8.88%  '

$ws.Range("D48").Value = '91.12'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("E49").Value = '  +5.83%  '

$ws.Range("E50").Value = '  +0.90%  '

$ws.Range("D51").Value = '8.22'
$ws.Range("E51").Value = '  +0.31%  '

# Restore the default cell style so formatting matches the original workbook.
$dataRange.Style = "Normal"
